$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.967.42'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.555.22'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").Value = '206.92'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '0.484'
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").Value = '21.67'
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = '0.0589'
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.775.86'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").Value = '1.554.15'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '3.71'
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").Value = '0.516'
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").Value = '26.957.04'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '61.82'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '214.81'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '0.0₃0688'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '7.27'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = '1.01'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '4.03'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = '153.26'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = '6.67'
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").Value = '14.89'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '0.104'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").Value = '0.0462'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").Value = '1.381.67'
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +2.81%  '
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  +6.03%  '
$ws.Range("D37").Value = '2.28'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("D39").Value = '0.523'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").Value = '0.809'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = '1.01'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").Value = '0.992'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '2.25'
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("D44").Value = '5.48'
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("D45").Value = '63.92'
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").Value = '1.689.94'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").Value = '86.09'
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").Value = '0.0956'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.45%  '
